$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the final existing row (row 75) with revised figures.
$ws.Range("B75").Value = -4.7
$ws.Range("C75").Value = -7.7
$ws.Range("D75").Value = -10.2

# Append a new row (row 76) for the next quarter.
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-04-2021"
$ws.Range("A76").ClearFormats()
$ws.Range("B76").Value = -4.3
$ws.Range("C76").Value = -7.2
$ws.Range("D76").Value = -9.5
$ws.Range("E76").Value = 34
